$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Inserting raw OOXML directly at the absolute end of the document's last
# paragraph clobbers that paragraph (Word COM quirk: the very last range in
# the body is special). Work around it by first appending a throwaway empty
# paragraph, then inserting our new content at the *start* of that
# paragraph's range -- the placeholder paragraph mark gets absorbed into the
# final inserted paragraph, leaving no stray empty paragraph behind.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$placeholder = $d.Paragraphs.Last
$r = $placeholder.Range
$r.Collapse(1)

$xml = @"
<w:p $ns>
  <w:pPr>
    <w:spacing w:line="320" w:lineRule="exact"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:spacing w:line="320" w:lineRule="exact"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve">Opinions of </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>d</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>efendant</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>'s</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve"> expert Dr. Doctor</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:spacing w:line="320" w:lineRule="exact"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>Here is what the doctor thought, and this is why we disagree with it.</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:line="320" w:lineRule="exact"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica Neue" w:hAnsi="Helvetica"/>
      <w:noProof/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica Neue" w:hAnsi="Helvetica"/>
      <w:noProof/>
    </w:rPr>
    <w:t>Number 1: disagreement 1</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:line="320" w:lineRule="exact"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica Neue" w:hAnsi="Helvetica"/>
      <w:noProof/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica Neue" w:hAnsi="Helvetica"/>
      <w:noProof/>
    </w:rPr>
    <w:t>Number 2: disagreement 2</w:t>
  </w:r>
</w:p>
"@

$r.InsertXML($xml)
Write-Output "inserted new paragraphs"
